$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateCell = $ws.Cells.Item(79, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "11/04/2025"
$dateCell.ClearFormats()

$ws.Cells.Item(79, 2).Value = 8821.99
